# Add missing historical daily quotes for 2019-11-18 .. 2019-11-28 into the
# middle of the price history table (between the existing 2019-11-15 and
# 2019-11-29 rows), shifting all subsequent rows down by 9 and extending the
# sheet's used range from A1:I156 to A1:I165.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 9 new rows that need to be inserted right before the current row 84
# (2019-11-29), in chronological order.
$newRows = @(
    @("1574035200", "2019-11-18", "0209", "ISTONE", 0.21,  0.21,  0.205, 0.205, "5369600"),
    @("1574121600", "2019-11-19", "0209", "ISTONE", 0.205, 0.22,  0.205, 0.22,  "22692900"),
    @("1574208000", "2019-11-20", "0209", "ISTONE", 0.225, 0.23,  0.22,  0.225, "17692000"),
    @("1574294400", "2019-11-21", "0209", "ISTONE", 0.225, 0.23,  0.22,  0.22,  "4626400"),
    @("1574380800", "2019-11-22", "0209", "ISTONE", 0.22,  0.225, 0.215, 0.22,  "7710200"),
    @("1574640000", "2019-11-25", "0209", "ISTONE", 0.22,  0.22,  0.21,  0.21,  "4889700"),
    @("1574726400", "2019-11-26", "0209", "ISTONE", 0.21,  0.22,  0.21,  0.215, "10732300"),
    @("1574812800", "2019-11-27", "0209", "ISTONE", 0.215, 0.225, 0.21,  0.225, "10991300"),
    @("1574899200", "2019-11-28", "0209", "ISTONE", 0.225, 0.23,  0.22,  0.225, "7400500")
)

$firstRow = 84
$lastRow = $firstRow + $newRows.Count - 1

# Push the existing rows 84:156 down by 9 rows (mirrors selecting rows
# 84-92, right-click, Insert in real Excel).
$ws.Range("A" + $firstRow + ":I" + $lastRow).EntireRow.Insert()

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $firstRow + $i
    $row = $newRows[$i]

    # Columns that must stay text (timestamp-looking / zero-padded id would
    # otherwise be auto-coerced by Excel into a date serial / plain number).
    $ws.Range("B" + $r + ":C" + $r).NumberFormat = "@"

    $ws.Range("A" + $r).Value = $row[0]
    $ws.Range("B" + $r).Value = $row[1]
    $ws.Range("C" + $r).Value = $row[2]
    $ws.Range("D" + $r).Value = $row[3]
    $ws.Range("E" + $r).Value = $row[4]
    $ws.Range("F" + $r).Value = $row[5]
    $ws.Range("G" + $r).Value = $row[6]
    $ws.Range("H" + $r).Value = $row[7]
    $ws.Range("I" + $r).Value = $row[8]

    # Drop the temporary Text number format again so the new rows end up
    # with the same (default/general) style as every other data row.
    $ws.Range("B" + $r + ":C" + $r).ClearFormats()
}
